$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-01 Sunday" "2026-02-02 Monday"

Replace-Text "40÷2=" "62÷3="
Replace-Text "34÷6=" "63÷6="
Replace-Text "47÷3=" "65÷7="
Replace-Text "53÷3=" "52÷4="
Replace-Text "84÷4=" "24÷9="
Replace-Text "32÷4=" "58÷2="
Replace-Text "85÷3=" "53÷6="
Replace-Text "65÷4=" "14÷3="
Replace-Text "42÷3=" "87÷3="
Replace-Text "46÷9=" "54÷9="
Replace-Text "71÷9=" "26÷7="
Replace-Text "45÷8=" "42÷9="
Replace-Text "94÷6=" "80÷2="
Replace-Text "11÷9=" "91÷5="
Replace-Text "20÷5=" "91÷9="
Replace-Text "69÷8=" "82÷7="
Replace-Text "98÷2=" "82÷6="
Replace-Text "68÷7=" "45÷9="
Replace-Text "60÷4=" "68÷4="
Replace-Text "17÷7=" "63÷3="
Replace-Text "70÷2=" "57÷8="
Replace-Text "49÷8=" "66÷6="
Replace-Text "27÷2=" "84÷2="
Replace-Text "60÷8=" "24÷9="
Replace-Text "58÷6=" "44÷8="
